$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.35"
$ws.Range("E2").Value = "'1.16%"
$ws.Range("G2").Value = "'3"

# Row 3
$ws.Range("D3").Value = "'29.38"
$ws.Range("E3").Value = "'-0.23%"
$ws.Range("G3").Value = "'3"

# Row 4
$ws.Range("D4").Value = "'5.167"
$ws.Range("E4").Value = "'0.97%"
$ws.Range("G4").Value = "'3"

# Row 5
$ws.Range("D5").Value = "'0.05760"
$ws.Range("E5").Value = "'1.88%"
$ws.Range("G5").Value = "'3"

# Row 6
$ws.Range("D6").Value = "'6.568"
$ws.Range("E6").Value = "'1.14%"
$ws.Range("G6").Value = "'3"

# Row 7
$ws.Range("D7").Value = "'0.8596"
$ws.Range("E7").Value = "'4.11%"
$ws.Range("G7").Value = "'3"

# Row 8
$ws.Range("D8").Value = "'0.8601"
$ws.Range("E8").Value = "'-0.55%"
$ws.Range("G8").Value = "'3"

# Row 9
$ws.Range("E9").Value = "'2.93%"
$ws.Range("G9").Value = "'3"

# Row 10
$ws.Range("D10").Value = "'0.07024"
$ws.Range("E10").Value = "'1.64%"
$ws.Range("G10").Value = "'3"

# Row 11
$ws.Range("D11").Value = "'0.03011"
$ws.Range("E11").Value = "'5.34%"
$ws.Range("G11").Value = "'3"

# Row 12
$ws.Range("D12").Value = "'0.09366"
$ws.Range("E12").Value = "'-0.25%"
$ws.Range("G12").Value = "'3"

# Row 13
$ws.Range("D13").Value = "'0.001526"
$ws.Range("E13").Value = "'0.88%"
$ws.Range("G13").Value = "'3"

# Row 14
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006031"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("G14").Value = "'3"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006120"
$ws.Range("E15").Value = "'0.31%"
$ws.Range("G15").Value = "'3"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.493"
$ws.Range("E16").Value = "'-0.79%"
$ws.Range("G16").Value = "'3"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.140"
$ws.Range("E17").Value = "'4.25%"
$ws.Range("G17").Value = "'3"

# Row 18
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.189"
$ws.Range("E18").Value = "'-1.27%"
$ws.Range("G18").Value = "'3"

# Row 19
$ws.Range("G19").Value = "'3"

# Row 20
$ws.Range("D20").Value = "'0.03291"
$ws.Range("E20").Value = "'1.38%"
$ws.Range("G20").Value = "'3"

# Row 21
$ws.Range("D21").Value = "'0.1284"
$ws.Range("E21").Value = "'-0.78%"
$ws.Range("G21").Value = "'3"

# Row 22
$ws.Range("D22").Value = "'3.176"
$ws.Range("E22").Value = "'-12.16%"
$ws.Range("G22").Value = "'3"

# Row 23
$ws.Range("E23").Value = "'-0.17%"
$ws.Range("G23").Value = "'3"

# Row 24
$ws.Range("E24").Value = "'1.96%"
$ws.Range("G24").Value = "'3"

# Row 25
$ws.Range("D25").Value = "'0.001227"
$ws.Range("E25").Value = "'1.42%"
$ws.Range("G25").Value = "'3"

# Row 26
$ws.Range("D26").Value = "'0.004131"
$ws.Range("E26").Value = "'-7.05%"
$ws.Range("G26").Value = "'3"

# Row 27
$ws.Range("D27").Value = "'0.0001211"
$ws.Range("E27").Value = "'2.62%"
$ws.Range("G27").Value = "'3"

# Row 28
$ws.Range("E28").Value = "'3.25%"
$ws.Range("G28").Value = "'3"

# Row 29
$ws.Range("G29").Value = "'3"

# Row 30
$ws.Range("G30").Value = "'3"

# Row 31
$ws.Range("G31").Value = "'3"

# Row 32
$ws.Range("G32").Value = "'3"

# Row 33
$ws.Range("G33").Value = "'3"

# Row 34
$ws.Range("G34").Value = "'3"

# Row 35
$ws.Range("G35").Value = "'3"

# Row 36
$ws.Range("G36").Value = "'3"

# Row 37
$ws.Range("G37").Value = "'3"

# Row 38
$ws.Range("G38").Value = "'3"

# Row 39
$ws.Range("G39").Value = "'3"

# Row 40
$ws.Range("D40").Value = "'0.03724"
$ws.Range("E40").Value = "'0.49%"
$ws.Range("G40").Value = "'3"

# Row 41
$ws.Range("D41").Value = "'0.005883"
$ws.Range("E41").Value = "'2.17%"
$ws.Range("G41").Value = "'3"

# Row 42
$ws.Range("D42").Value = "'0.1070"
$ws.Range("E42").Value = "'1.75%"
$ws.Range("G42").Value = "'3"

# Row 43
$ws.Range("D43").Value = "'0.002441"
$ws.Range("E43").Value = "'5.64%"
$ws.Range("G43").Value = "'3"

# Row 44
$ws.Range("D44").Value = "'0.008399"
$ws.Range("E44").Value = "'-13.46%"
$ws.Range("G44").Value = "'3"

# Row 45
$ws.Range("D45").Value = "'0.00005274"
$ws.Range("E45").Value = "'3.45%"
$ws.Range("G45").Value = "'3"

# Row 46
$ws.Range("E46").Value = "'0.04%"
$ws.Range("G46").Value = "'3"

# Row 47
$ws.Range("D47").Value = "'0.05801"
$ws.Range("E47").Value = "'-44.74%"
$ws.Range("G47").Value = "'3"

# Row 48
$ws.Range("D48").Value = "'0.002447"
$ws.Range("E48").Value = "'-3.83%"
$ws.Range("G48").Value = "'3"

# Row 49
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("G49").Value = "'3"

# Row 50
$ws.Range("E50").Value = "'0.04%"
$ws.Range("G50").Value = "'3"

# Row 51
$ws.Range("G51").Value = "'3"
